$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, pushing existing rows 43-62 down to 44-63.
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the latest weekly price entry.
$ws.Cells.Item(43, 1).Value = 5
$ws.Cells.Item(43, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(43, 3).Value = "Maule"
$ws.Cells.Item(43, 4).Value = 44460
$ws.Cells.Item(43, 5).Value = 7
$ws.Cells.Item(43, 6).Value = 100112001
$ws.Cells.Item(43, 7).Value = "Berenjena"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 200
$ws.Cells.Item(43, 11).Value = 7000
$ws.Cells.Item(43, 12).Value = 7000
$ws.Cells.Item(43, 13).Value = 7000
$ws.Cells.Item(43, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(43, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(43, 16).Value = 117
$ws.Cells.Item(43, 17).Value = 60
$ws.Cells.Item(43, 18).Value = "Hortaliza"
